$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.660188000000001
$ws.Range("H2").Value = 19.980564
$ws.Range("I2").Value = 0.1500148400131262
$ws.Range("J2").Value = 0.1500148400131261
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.664391
$ws.Range("N2").Value = 4.993173000000001
$ws.Range("O2").Value = 0.3990511495040125
$ws.Range("P2").Value = 0.3990511495040125
$ws.Range("Q2").Value = 11.085156965508
$ws.Range("R2").Value = 99.76641268957202
$ws.Range("S2").Value = 0.05986359434989854
$ws.Range("T2").Value = 0.05986359434989853

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.660188000000001
$ws.Range("H3").Value = 19.980564
$ws.Range("I3").Value = 0.1500148400131262
$ws.Range("J3").Value = 0.1500148400131261
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.437958
$ws.Range("N3").Value = 1.313874
$ws.Range("O3").Value = 0.1050039584054939
$ws.Range("P3").Value = 0.1050039584054938
$ws.Range("Q3").Value = 2.916882616104
$ws.Range("R3").Value = 26.251943544936
$ws.Range("S3").Value = 0.01575215202094512
$ws.Range("T3").Value = 0.01575215202094511

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.660188000000001
$ws.Range("H4").Value = 19.980564
$ws.Range("I4").Value = 0.1500148400131262
$ws.Range("J4").Value = 0.1500148400131261
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.439215333333333
$ws.Range("N4").Value = 4.317646
$ws.Range("O4").Value = 0.3450634695516061
$ws.Range("P4").Value = 0.3450634695516061
$ws.Range("Q4").Value = 9.585444692482668
$ws.Range("R4").Value = 86.26900223234401
$ws.Range("S4").Value = 0.05176464117915842
$ws.Range("T4").Value = 0.05176464117915841

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.660188000000001
$ws.Range("H5").Value = 19.980564
$ws.Range("I5").Value = 0.1500148400131262
$ws.Range("J5").Value = 0.1500148400131261
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.629307
$ws.Range("N5").Value = 1.887921
$ws.Range("O5").Value = 0.1508814225388875
$ws.Range("P5").Value = 0.1508814225388875
$ws.Range("Q5").Value = 4.191302929716
$ws.Range("R5").Value = 37.721726367444
$ws.Range("S5").Value = 0.0226344524631241
$ws.Range("T5").Value = 0.0226344524631241

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 16.410331
$ws.Range("H6").Value = 49.230993
$ws.Range("I6").Value = 0.3696281815959916
$ws.Range("J6").Value = 0.3696281815959916
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.664391
$ws.Range("N6").Value = 4.993173000000001
$ws.Range("O6").Value = 0.3990511495040125
$ws.Range("P6").Value = 0.3990511495040125
$ws.Range("Q6").Value = 27.313207223421
$ws.Range("R6").Value = 245.818865010789
$ws.Range("S6").Value = 0.1475005507549584
$ws.Range("T6").Value = 0.1475005507549584

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 16.410331
$ws.Range("H7").Value = 49.230993
$ws.Range("I7").Value = 0.3696281815959916
$ws.Range("J7").Value = 0.3696281815959916
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.437958
$ws.Range("N7").Value = 1.313874
$ws.Range("O7").Value = 0.1050039584054939
$ws.Range("P7").Value = 0.1050039584054938
$ws.Range("Q7").Value = 7.187035744098
$ws.Range("R7").Value = 64.68332169688199
$ws.Range("S7").Value = 0.03881242220580384
$ws.Range("T7").Value = 0.03881242220580383

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 16.410331
$ws.Range("H8").Value = 49.230993
$ws.Range("I8").Value = 0.3696281815959916
$ws.Range("J8").Value = 0.3696281815959916
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.439215333333333
$ws.Range("N8").Value = 4.317646
$ws.Range("O8").Value = 0.3450634695516061
$ws.Range("P8").Value = 0.3450634695516061
$ws.Range("Q8").Value = 23.61800000027533
$ws.Range("R8").Value = 212.562000002478
$ws.Range("S8").Value = 0.127545182785564
$ws.Range("T8").Value = 0.127545182785564

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 16.410331
$ws.Range("H9").Value = 49.230993
$ws.Range("I9").Value = 0.3696281815959916
$ws.Range("J9").Value = 0.3696281815959916
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.629307
$ws.Range("N9").Value = 1.887921
$ws.Range("O9").Value = 0.1508814225388875
$ws.Range("P9").Value = 0.1508814225388875
$ws.Range("Q9").Value = 10.327136170617
$ws.Range("R9").Value = 92.94422553555299
$ws.Range("S9").Value = 0.05577002584966547
$ws.Range("T9").Value = 0.05577002584966547

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.478895
$ws.Range("H10").Value = 43.436685
$ws.Range("I10").Value = 0.3261242951387937
$ws.Range("J10").Value = 0.3261242951387937
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.664391
$ws.Range("N10").Value = 4.993173000000001
$ws.Range("O10").Value = 0.3990511495040125
$ws.Range("P10").Value = 0.3990511495040125
$ws.Range("Q10").Value = 24.098542527945
$ws.Range("R10").Value = 216.886882751505
$ws.Range("S10").Value = 0.1301402748563215
$ws.Range("T10").Value = 0.1301402748563215

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 14.478895
$ws.Range("H11").Value = 43.436685
$ws.Range("I11").Value = 0.3261242951387937
$ws.Range("J11").Value = 0.3261242951387937
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.437958
$ws.Range("N11").Value = 1.313874
$ws.Range("O11").Value = 0.1050039584054939
$ws.Range("P11").Value = 0.1050039584054938
$ws.Range("Q11").Value = 6.34114789641
$ws.Range("R11").Value = 57.07033106769
$ws.Range("S11").Value = 0.0342443419217749
$ws.Range("T11").Value = 0.03424434192177488

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 14.478895
$ws.Range("H12").Value = 43.436685
$ws.Range("I12").Value = 0.3261242951387937
$ws.Range("J12").Value = 0.3261242951387937
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.439215333333333
$ws.Range("N12").Value = 4.317646
$ws.Range("O12").Value = 0.3450634695516061
$ws.Range("P12").Value = 0.3450634695516061
$ws.Range("Q12").Value = 20.83824769372333
$ws.Range("R12").Value = 187.54422924351
$ws.Range("S12").Value = 0.1125335807856641
$ws.Range("T12").Value = 0.1125335807856641

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 14.478895
$ws.Range("H13").Value = 43.436685
$ws.Range("I13").Value = 0.3261242951387937
$ws.Range("J13").Value = 0.3261242951387937
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.629307
$ws.Range("N13").Value = 1.887921
$ws.Range("O13").Value = 0.1508814225388875
$ws.Range("P13").Value = 0.1508814225388875
$ws.Range("Q13").Value = 9.111669975764999
$ws.Range("R13").Value = 82.005029781885
$ws.Range("S13").Value = 0.0492060975750332
$ws.Range("T13").Value = 0.04920609757503319

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.847447
$ws.Range("H14").Value = 20.542341
$ws.Range("I14").Value = 0.1542326832520885
$ws.Range("J14").Value = 0.1542326832520885
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.664391
$ws.Range("N14").Value = 4.993173000000001
$ws.Range("O14").Value = 0.3990511495040125
$ws.Range("P14").Value = 0.3990511495040125
$ws.Range("Q14").Value = 11.396829159777
$ws.Range("R14").Value = 102.571462437993
$ws.Range("S14").Value = 0.06154672954283417
$ws.Range("T14").Value = 0.06154672954283417

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.847447
$ws.Range("H15").Value = 20.542341
$ws.Range("I15").Value = 0.1542326832520885
$ws.Range("J15").Value = 0.1542326832520885
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.437958
$ws.Range("N15").Value = 1.313874
$ws.Range("O15").Value = 0.1050039584054939
$ws.Range("P15").Value = 0.1050039584054938
$ws.Range("Q15").Value = 2.998894193226
$ws.Range("R15").Value = 26.990047739034
$ws.Range("S15").Value = 0.01619504225697001
$ws.Range("T15").Value = 0.01619504225697

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.847447
$ws.Range("H16").Value = 20.542341
$ws.Range("I16").Value = 0.1542326832520885
$ws.Range("J16").Value = 0.1542326832520885
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.439215333333333
$ws.Range("N16").Value = 4.317646
$ws.Range("O16").Value = 0.3450634695516061
$ws.Range("P16").Value = 0.3450634695516061
$ws.Range("Q16").Value = 9.854950716587332
$ws.Range("R16").Value = 88.694556449286
$ws.Range("S16").Value = 0.05322006480121953
$ws.Range("T16").Value = 0.05322006480121953

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.847447
$ws.Range("H17").Value = 20.542341
$ws.Range("I17").Value = 0.1542326832520885
$ws.Range("J17").Value = 0.1542326832520885
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.629307
$ws.Range("N17").Value = 1.887921
$ws.Range("O17").Value = 0.1508814225388875
$ws.Range("P17").Value = 0.1508814225388875
$ws.Range("Q17").Value = 4.309146329229
$ws.Range("R17").Value = 38.782316963061
$ws.Range("S17").Value = 0.02327084665106476
$ws.Range("T17").Value = 0.02327084665106476

